# Auto-generated edit script: updates numeric market/profit data cells
# across multiple worksheets (ARM, BSM, CRP, CUL, GSM, LTW) per the commit diff.
$wb = $excel.ActiveWorkbook

# ----- Sheet: ARM -----
$ws = $wb.Worksheets.Item("ARM")

$ws.Range("H5").Value = 294.76923
$ws.Range("I5").Value = 106.666664
$ws.Range("J5").Value = 456
$ws.Range("K5").Value = 106.666664
$ws.Range("L5").Value = 456
$ws.Range("M5").Value = 5.333336000000003
$ws.Range("N5").Value = -680
$ws.Range("H9").Value = 44000
$ws.Range("I9").Value = 44000
$ws.Range("K9").Value = 44000
$ws.Range("M9").Value = -43830
$ws.Range("H20").Value = 44000
$ws.Range("I20").Value = 44000
$ws.Range("K20").Value = 44000
$ws.Range("M20").Value = -43730
$ws.Range("H32").Value = 6450013.5
$ws.Range("I32").Value = 1549110.6
$ws.Range("J32").Value = 31260836
$ws.Range("K32").Value = 1549110.6
$ws.Range("L32").Value = 31260836
$ws.Range("M32").Value = -1548823.6
$ws.Range("N32").Value = -31261410
$ws.Range("H37").Value = 11077.733
$ws.Range("I37").Value = 5780
$ws.Range("J37").Value = 13726.6
$ws.Range("K37").Value = 5780
$ws.Range("L37").Value = 13726.6
$ws.Range("M37").Value = -5507
$ws.Range("N37").Value = -14272.6
$ws.Range("H44").Value = 17714.285
$ws.Range("I44").Value = 2000
$ws.Range("J44").Value = 20333.334
$ws.Range("K44").Value = 2000
$ws.Range("L44").Value = 20333.334
$ws.Range("M44").Value = -1512
$ws.Range("N44").Value = -21309.334
$ws.Range("H55").Value = 18456.555
$ws.Range("I55").Value = 6333.3335
$ws.Range("J55").Value = 24518.166
$ws.Range("K55").Value = 6333.3335
$ws.Range("L55").Value = 24518.166
$ws.Range("M55").Value = -6018.3335
$ws.Range("N55").Value = -25148.166
$ws.Range("H63").Value = 1735.909
$ws.Range("I63").Value = 1709.5
$ws.Range("J63").Value = 2000
$ws.Range("K63").Value = 1709.5
$ws.Range("L63").Value = 2000
$ws.Range("M63").Value = -1023.5
$ws.Range("N63").Value = -3372
$ws.Range("H66").Value = 1735.909
$ws.Range("I66").Value = 1709.5
$ws.Range("J66").Value = 2000
$ws.Range("K66").Value = 8547.5
$ws.Range("L66").Value = 10000
$ws.Range("M66").Value = -5115.5
$ws.Range("N66").Value = -16864
$ws.Range("H80").Value = 38226
$ws.Range("J80").Value = 38226
$ws.Range("L80").Value = 38226
$ws.Range("N80").Value = -40222
$ws.Range("H83").Value = 38226
$ws.Range("J83").Value = 38226
$ws.Range("L83").Value = 114678
$ws.Range("N83").Value = -124662
$ws.Range("H122").Value = 3720.3635
$ws.Range("I122").Value = 3378
$ws.Range("J122").Value = 4633.3335
$ws.Range("K122").Value = 10134
$ws.Range("L122").Value = 13900.0005
$ws.Range("M122").Value = -7684
$ws.Range("N122").Value = -18800.0005

# ----- Sheet: BSM -----
$ws = $wb.Worksheets.Item("BSM")

$ws.Range("H4").Value = 294.76923
$ws.Range("I4").Value = 106.666664
$ws.Range("J4").Value = 456
$ws.Range("K4").Value = 106.666664
$ws.Range("L4").Value = 456
$ws.Range("M4").Value = 8.333336000000003
$ws.Range("N4").Value = -686
$ws.Range("H15").Value = 13500
$ws.Range("I15").Value = 13500
$ws.Range("K15").Value = 13500
$ws.Range("M15").Value = -13273
$ws.Range("H19").Value = 0
$ws.Range("I19").Value = 0
$ws.Range("K19").Value = 0
$ws.Range("H35").Value = 0
$ws.Range("J35").Value = 0
$ws.Range("L35").Value = 0
$ws.Range("H82").Value = 20000
$ws.Range("J82").Value = 25000
$ws.Range("L82").Value = 25000
$ws.Range("N82").Value = -25766
$ws.Range("H85").Value = 20000
$ws.Range("J85").Value = 25000
$ws.Range("L85").Value = 25000
$ws.Range("N85").Value = -27652

# Cell(s) removed entirely in the source edit (clear so no cell remains)
$ws.Range("M19").ClearContents()
$ws.Range("N35").ClearContents()

# ----- Sheet: CRP -----
$ws = $wb.Worksheets.Item("CRP")

$ws.Range("H16").Value = 635.7143
$ws.Range("I16").Value = 540.0909
$ws.Range("J16").Value = 986.3333
$ws.Range("K16").Value = 540.0909
$ws.Range("L16").Value = 986.3333
$ws.Range("M16").Value = -253.0909
$ws.Range("N16").Value = -1560.3333
$ws.Range("H31").Value = 2217736
$ws.Range("I31").Value = 1017569.9
$ws.Range("K31").Value = 1017569.9
$ws.Range("M31").Value = -1017274.9
$ws.Range("H34").Value = 2217736
$ws.Range("I34").Value = 1017569.9
$ws.Range("K34").Value = 1017569.9
$ws.Range("M34").Value = -1017367.9
$ws.Range("H38").Value = 29999
$ws.Range("I38").Value = 0
$ws.Range("J38").Value = 29999
$ws.Range("K38").Value = 0
$ws.Range("L38").Value = 29999
$ws.Range("N38").Value = -30753
$ws.Range("H41").Value = 7637.857
$ws.Range("I41").Value = 5583.3335
$ws.Range("J41").Value = 19965
$ws.Range("K41").Value = 5583.3335
$ws.Range("L41").Value = 19965
$ws.Range("M41").Value = -5155.3335
$ws.Range("N41").Value = -20821
$ws.Range("H46").Value = 29999
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 29999
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 29999
$ws.Range("M46").Value = 29999
$ws.Range("N46").Value = -30421
$ws.Range("H51").Value = 8199.111000000001
$ws.Range("I51").Value = 0
$ws.Range("J51").Value = 8199.111000000001
$ws.Range("K51").Value = 0
$ws.Range("L51").Value = 8199.111000000001
$ws.Range("N51").Value = -9671.111000000001
$ws.Range("H60").Value = 10997.5
$ws.Range("I60").Value = 2000
$ws.Range("J60").Value = 19995
$ws.Range("K60").Value = 2000
$ws.Range("L60").Value = 19995
$ws.Range("M60").Value = -1489
$ws.Range("N60").Value = -21017
$ws.Range("H61").Value = 8199.111000000001
$ws.Range("I61").Value = 0
$ws.Range("J61").Value = 8199.111000000001
$ws.Range("K61").Value = 0
$ws.Range("L61").Value = 8199.111000000001
$ws.Range("N61").Value = -8895.111000000001
$ws.Range("H68").Value = 18406.5
$ws.Range("J68").Value = 17340.555
$ws.Range("L68").Value = 17340.555
$ws.Range("N68").Value = -18838.555
$ws.Range("H71").Value = 18406.5
$ws.Range("J71").Value = 17340.555
$ws.Range("L71").Value = 52021.665
$ws.Range("N71").Value = -59509.665
$ws.Range("H74").Value = 17982.908
$ws.Range("J74").Value = 17982.908
$ws.Range("L74").Value = 17982.908
$ws.Range("N74").Value = -19730.908
$ws.Range("H77").Value = 17982.908
$ws.Range("J77").Value = 17982.908
$ws.Range("L77").Value = 53948.724
$ws.Range("N77").Value = -62684.724
$ws.Range("H113").Value = 635.7143
$ws.Range("I113").Value = 540.0909
$ws.Range("J113").Value = 986.3333
$ws.Range("K113").Value = 540.0909
$ws.Range("L113").Value = 986.3333
$ws.Range("M113").Value = 1629.9091
$ws.Range("N113").Value = -5326.3333
$ws.Range("H132").Value = 1906.1111
$ws.Range("I132").Value = 1545.0344
$ws.Range("J132").Value = 3402
$ws.Range("K132").Value = 4635.1032
$ws.Range("L132").Value = 10206
$ws.Range("M132").Value = -2105.1032
$ws.Range("N132").Value = -15266

# Cell(s) removed entirely in the source edit (clear so no cell remains)
$ws.Range("M38").ClearContents()
$ws.Range("M51").ClearContents()
$ws.Range("M61").ClearContents()

# ----- Sheet: CUL -----
$ws = $wb.Worksheets.Item("CUL")

$ws.Range("H2").Value = 53.214287
$ws.Range("I2").Value = 38.5
$ws.Range("J2").Value = 72.833336
$ws.Range("K2").Value = 231
$ws.Range("L2").Value = 437.000016
$ws.Range("M2").Value = -118
$ws.Range("N2").Value = -663.000016
$ws.Range("H17").Value = 250
$ws.Range("I17").Value = 250
$ws.Range("J17").Value = 0
$ws.Range("K17").Value = 750
$ws.Range("L17").Value = 0
$ws.Range("M17").Value = -581
$ws.Range("H34").Value = 1228.4
$ws.Range("I34").Value = 677.8
$ws.Range("J34").Value = 1779
$ws.Range("K34").Value = 2033.4
$ws.Range("L34").Value = 5337
$ws.Range("M34").Value = -1949.4
$ws.Range("N34").Value = -5505
$ws.Range("H39").Value = 2537.0833
$ws.Range("J39").Value = 2537.0833
$ws.Range("L39").Value = 7611.249899999999
$ws.Range("N39").Value = -8199.249899999999
$ws.Range("H55").Value = 3191.8462
$ws.Range("J55").Value = 3415.8333
$ws.Range("L55").Value = 10247.4999
$ws.Range("N55").Value = -10601.4999
$ws.Range("H132").Value = 2196
$ws.Range("I132").Value = 3156.4
$ws.Range("J132").Value = 1715.8
$ws.Range("K132").Value = 28407.6
$ws.Range("L132").Value = 15442.2
$ws.Range("M132").Value = -25877.6
$ws.Range("N132").Value = -20502.2

# Cell(s) removed entirely in the source edit (clear so no cell remains)
$ws.Range("N17").ClearContents()

# ----- Sheet: GSM -----
$ws = $wb.Worksheets.Item("GSM")

$ws.Range("H113").Value = 25241.455
$ws.Range("I113").Value = 2765.6
$ws.Range("J113").Value = 250000
$ws.Range("K113").Value = 2765.6
$ws.Range("L113").Value = 250000
$ws.Range("M113").Value = -595.5999999999999
$ws.Range("N113").Value = -254340

# ----- Sheet: LTW -----
$ws = $wb.Worksheets.Item("LTW")

$ws.Range("H22").Value = 2559.1765
$ws.Range("I22").Value = 1712.5
$ws.Range("J22").Value = 3311.7778
$ws.Range("K22").Value = 1712.5
$ws.Range("L22").Value = 3311.7778
$ws.Range("M22").Value = -1417.5
$ws.Range("N22").Value = -3901.7778
$ws.Range("H27").Value = 2559.1765
$ws.Range("I27").Value = 1712.5
$ws.Range("J27").Value = 3311.7778
$ws.Range("K27").Value = 1712.5
$ws.Range("L27").Value = 3311.7778
$ws.Range("M27").Value = -1605.5
$ws.Range("N27").Value = -3525.7778
$ws.Range("H46").Value = 536.5862
$ws.Range("I46").Value = 515
$ws.Range("J46").Value = 551.82355
$ws.Range("K46").Value = 515
$ws.Range("L46").Value = 551.82355
$ws.Range("M46").Value = -327
$ws.Range("N46").Value = -927.82355
$ws.Range("H100").Value = 2201.25
$ws.Range("I100").Value = 1242.8572
$ws.Range("J100").Value = 2946.6667
$ws.Range("K100").Value = 1242.8572
$ws.Range("L100").Value = 2946.6667
$ws.Range("M100").Value = -701.8571999999999
$ws.Range("N100").Value = -4028.6667

